# The sheet is renamed ("Test1" -> "ValidLogin") and, per the diff, its
# sheetId advances from 1 to 2. A plain in-place rename keeps sheetId==1,
# so instead we add a brand new sheet (which is allocated the next
# sheetId, i.e. 2), populate it, and then remove the old "Test1" sheet.
$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "ValidLogin"

# Shared-string/cell content updates:
#   username -> UserName, password -> Password, bhanu -> admin, akshara -> manager
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "manager"

# Remove the original sheet now that its data/name live on the new sheet.
$oldSheet = $wb.Worksheets.Item("Test1")
$oldSheet.Delete()

# Restore the view state: selection moves from B3 to B1, zoom 235% -> 160%.
$newSheet.Range("B1").Select()
$excel.ActiveWindow.Zoom = 160
